# Convert the "Metal" column (E) on the "Elements Table" sheet from a
# free-text "yes"/blank string column into a proper TRUE/FALSE boolean
# column, and leave the selection on E3 (matching the author's edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Elements Table")

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

# Data starts on row 3 (row 1 = field-name header, row 2 = type header).
for ($r = 3; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $txt = $cell.Value2
    if ($txt -eq "yes") {
        $cell.Value = $true
    } else {
        $cell.Value = $false
    }
}

[void]$ws.Range("E3").Select()
